$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Project Coordinator" + ": Kyle Kalmon" -> single merged run
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Project Coordinator: Kyle Kalmon", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Project Coordinator: Kyle Kalmon", 2)

# ---------------------------------------------------------------------------
# 2) "QA Czar: Matt " / proofErr / "McCullar" / proofErr -> single merged run
#    with no leftover proofErr markers. Build a clean replacement paragraph
#    right after it, then delete the old (proofErr-laden) paragraph outright.
# ---------------------------------------------------------------------------
$qaOld = $d.Paragraphs(4)
$qaOld.Range.InsertParagraphAfter()
$qaNew = $d.Paragraphs(5)
$qaNew.Range.Text = "QA Czar: Matt McCullar"

$qaOld = $d.Paragraphs(4)
$qaOldRange = $qaOld.Range
$afterOld = $d.Paragraphs(5).Range
$killRange = $d.Range($qaOldRange.Start, $afterOld.Start)
$killRange.Delete()

# QA Czar paragraph is now paragraph 4, clean, no proofErr.
$qaCzarPara = $d.Paragraphs(4)

# ---------------------------------------------------------------------------
# 3) Insert four new ListParagraph bullet items after the QA Czar line:
#      Code Skeleton - Matt McCullar
#      UML(Models)/Class diagram - Kyle Kalmon   (3 runs + _GoBack bookmark)
#      SiteMap - Terrell Martin
#      Wireframe - Alex
#    All runs get explicit black font color (w:color val="000000").
# ---------------------------------------------------------------------------

# -- Code Skeleton - Matt McCullar -----------------------------------------
$qaCzarPara.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs(5)
$p1.Range.Text = "Code Skeleton - Matt McCullar"
$p1r = $d.Paragraphs(5).Range
$p1text = $d.Range($p1r.Start, $p1r.End - 1)
$p1text.Font.Color = 0

# -- UML(Models)/Class diagram - Kyle Kalmon -------------------------------
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(6)
$p2.Range.Text = "UML(Models)/Class diagram - Kyle Kalmon"
$p2r = $d.Paragraphs(6).Range
$p2full = $d.Range($p2r.Start, $p2r.End - 1)

$len1 = "UML(Models)".Length
$len2 = "/Class diagram".Length

$run1 = $d.Range($p2full.Start, $p2full.Start + $len1)
$run2 = $d.Range($p2full.Start + $len1, $p2full.Start + $len1 + $len2)
$run3 = $d.Range($p2full.Start + $len1 + $len2, $p2full.End)
$run1.Font.Color = 0
$run2.Font.Color = 0
$run3.Font.Color = 0

$bmPos = $p2full.Start + $len1 + $len2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# -- SiteMap - Terrell Martin ------------------------------------------------
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs(7)
$p3.Range.Text = "SiteMap - Terrell Martin"
$p3r = $d.Paragraphs(7).Range
$p3text = $d.Range($p3r.Start, $p3r.End - 1)
$p3text.Font.Color = 0

# -- Wireframe - Alex ---------------------------------------------------------
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs(8)
$p4.Range.Text = "Wireframe - Alex"
$p4r = $d.Paragraphs(8).Range
$p4text = $d.Range($p4r.Start, $p4r.End - 1)
$p4text.Font.Color = 0

# ---------------------------------------------------------------------------
# 4) "Iteration 2" + bookmark(_GoBack) + ":" -> single merged run, bookmark
#    gone (it moved earlier in the doc, handled above).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Iteration 2:", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Iteration 2:", 2)

Write-Output "done"
